$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue $ws 'D2' '43.645.59'
Set-TextValue $ws 'E2' '  -0.63%  '
Set-TextValue $ws 'D3' '2.333.84'
Set-TextValue $ws 'E3' '  -1.14%  '
Set-TextValue $ws 'E4' '  +0.07%  '
Set-TextValue $ws 'D5' '238.34'
Set-TextValue $ws 'E5' '  -1.32%  '
Set-TextValue $ws 'E6' '  -4.56%  '
Set-TextValue $ws 'D7' '71.54'
Set-TextValue $ws 'E7' '  -6.48%  '
Set-TextValue $ws 'E8' '  +0.03%  '
Set-TextValue $ws 'D9' '0.579'
Set-TextValue $ws 'E9' '  -9.06%  '
Set-TextValue $ws 'D10' '0.0980'
Set-TextValue $ws 'E10' '  -4.94%  '
Set-TextValue $ws 'D11' '57.94'
Set-TextValue $ws 'E11' '  +1.05%  '
Set-TextValue $ws 'D12' '32.23'
Set-TextValue $ws 'E12' '  -3.71%  '
Set-TextValue $ws 'D13' '0.108'
Set-TextValue $ws 'E13' '  -0.61%  '
Set-TextValue $ws 'D14' '7.07'
Set-TextValue $ws 'E14' '  -6.85%  '
Set-TextValue $ws 'D15' '2.680.41'
Set-TextValue $ws 'E15' '  -1.14%  '
Set-TextValue $ws 'D16' '15.94'
Set-TextValue $ws 'E16' '  -5.83%  '
Set-TextValue $ws 'E17' '  -3.87%  '
Set-TextValue $ws 'D18' '2.338.76'
Set-TextValue $ws 'E18' '  -0.64%  '
Set-TextValue $ws 'D19' '43.597.59'
Set-TextValue $ws 'E19' '  -0.67%  '
Set-TextValue $ws 'D20' '0.0000100'
Set-TextValue $ws 'E20' '  -4.29%  '
Set-TextValue $ws 'D21' '77.82'
Set-TextValue $ws 'E21' '  +0.37%  '
Set-TextValue $ws 'E22' '  -1.92%  '
Set-TextValue $ws 'D23' '249.82'
Set-TextValue $ws 'E23' '  -2.80%  '
Set-TextValue $ws 'B24' 'Dai'
Set-TextValue $ws 'C24' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws 'D24' '1.00'
Set-TextValue $ws 'E24' '  -0.05%  '
Set-TextValue $ws 'B25' 'ImmutableX'
Set-TextValue $ws 'C25' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws 'D25' '1.90'
Set-TextValue $ws 'E25' '  +6.80%  '
Set-TextValue $ws 'E26' '  +2.97%  '
Set-TextValue $ws 'E27' '  -2.05%  '
Set-TextValue $ws 'D28' '10.24'
Set-TextValue $ws 'E28' '  -8.63%  '
Set-TextValue $ws 'E29' '  -1.60%  '
Set-TextValue $ws 'D30' '175.56'
Set-TextValue $ws 'E30' '  +0.42%  '
Set-TextValue $ws 'D31' '21.99'
Set-TextValue $ws 'E31' '  -5.16%  '
Set-TextValue $ws 'E32' '  -2.92%  '
Set-TextValue $ws 'E33' '  -1.43%  '
Set-TextValue $ws 'D34' '0.0728'
Set-TextValue $ws 'E34' '  -3.31%  '
Set-TextValue $ws 'E35' '  -5.73%  '
Set-TextValue $ws 'D36' '5.30'
Set-TextValue $ws 'E36' '  -1.07%  '
Set-TextValue $ws 'E37' '  -2.88%  '
Set-TextValue $ws 'E38' '  -1.86%  '
Set-TextValue $ws 'E39' '  -3.93%  '
Set-TextValue $ws 'D40' '5.60'
Set-TextValue $ws 'E40' '  +24.18%  '
Set-TextValue $ws 'E41' '  -3.07%  '
Set-TextValue $ws 'D42' '64.88'
Set-TextValue $ws 'E42' '  +18.35%  '
Set-TextValue $ws 'E43' '  +2.46%  '
Set-TextValue $ws 'E44' '  +3.20%  '
Set-TextValue $ws 'D45' '18.64'
Set-TextValue $ws 'E45' '  -3.35%  '
Set-TextValue $ws 'E46' '  -4.26%  '
Set-TextValue $ws 'E47' '  +0.03%  '
Set-TextValue $ws 'B48' 'HuobiToken'
Set-TextValue $ws 'C48' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws 'D48' '2.94'
Set-TextValue $ws 'E48' '  +4.50%  '
Set-TextValue $ws 'B49' 'TrustWalletToken'
Set-TextValue $ws 'C49' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws 'D49' '1.21'
Set-TextValue $ws 'E49' '  -4.18%  '
Set-TextValue $ws 'E50' '  -4.27%  '
Set-TextValue $ws 'D51' '97.33'
Set-TextValue $ws 'E51' '  -4.69%  '
